$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.781999999999999
$ws.Range("B7").Value = 6.032
$ws.Range("D7").Value = -7.567
$ws.Range("D15").Value = -8.271000000000001
$ws.Range("B16").Value = 5.315
$ws.Range("D21").Value = -8.1
$ws.Range("D22").Value = -7.948
$ws.Range("D23").Value = -7.869999999999999
$ws.Range("B28").Value = 6.034000000000001
$ws.Range("B29").Value = 5.615
$ws.Range("B32").Value = 6.572
$ws.Range("D34").Value = -7.917999999999999
$ws.Range("B40").Value = 9.223000000000001
$ws.Range("D43").Value = -7.712000000000001
$ws.Range("D45").Value = -7.532999999999999
$ws.Range("D50").Value = -8.096
$ws.Range("D51").Value = -8.180000000000001
$ws.Range("B52").Value = 5.359
$ws.Range("B57").Value = 5.090999999999999
$ws.Range("B66").Value = 5.773
$ws.Range("D66").Value = -7.561
$ws.Range("D67").Value = -7.157999999999999
$ws.Range("D79").Value = -7.446000000000001
$ws.Range("D84").Value = -8.132000000000001
$ws.Range("D92").Value = -6.602999999999999
$ws.Range("D97").Value = -8.416
$ws.Range("B100").Value = 5.931
